$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# from 2023-10-04 (45203) to 2023-10-06 (45205) for every data row (2..271).
$ws.Range("C2:C271").Value = 45205
